$d = $word.ActiveDocument

# Locate the three consecutive paragraphs to remove:
#   1) the blank paragraph right after "LOB1038: Física Experimental I (Requisito fraco)"
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) the "© 2020 . Contact: ..." footer paragraph
# Build a single Range spanning from the start of paragraph #1 through the
# end of paragraph #3 (inclusive of its trailing paragraph mark) and delete it.

$startPara = $null
$endPara = $null

$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($t -like "*LOB1038: Física Experimental I (Requisito fraco)*") {
        $startPara = $d.Paragraphs.Item($i + 1)
    }
    if ($t -like "*Contact: luizeleno@usp.br*") {
        $endPara = $p
    }
}

$range = $d.Range($startPara.Range.Start, $endPara.Range.End)
$range.Delete()
